# WS_holdings.xlsx update:
#  - bump the "as of" date in the confidential disclaimer (A16) from
#    2021-05-18 to 2021-05-19
#  - refresh the Weight (D) / Percent Change (E) values for rows 2-13
#
# The sheet ships protected (no known password), so it has to be
# unprotected before the locked cells can be written, and re-protected
# afterwards to restore the original behaviour.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Unprotect()

$ws.Range("A16").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-19 for illustrative purposes only and are subject to change."

$ws.Range("D2").Value = 0.03103685307233568
$ws.Range("E2").Value = -0.004123711340206171

$ws.Range("D3").Value = 0.02378557332194393
$ws.Range("E3").Value = -0.001183992422448488

$ws.Range("D4").Value = 0.05263430368588246
$ws.Range("E4").Value = -0.003715745471435161

$ws.Range("D5").Value = 0.1387954551664335
$ws.Range("E5").Value = -0.002937336814621494

$ws.Range("D6").Value = 0.03247167388331323
$ws.Range("E6").Value = -0.02484472049689446

$ws.Range("D7").Value = 0.1175396627949394
$ws.Range("E7").Value = -0.01048815227243305

$ws.Range("D8").Value = 0.102624500934497
$ws.Range("E8").Value = -0.007195571955719626

$ws.Range("D9").Value = 0.03000253097933082
$ws.Range("E9").Value = -0.01588918313302101

$ws.Range("D10").Value = 0.1279151201568057
$ws.Range("E10").Value = -0.005433764287052667

$ws.Range("D11").Value = 0.2416042670010628
$ws.Range("E11").Value = 0.002921496560173376

$ws.Range("D12").Value = 0.1015900590034557
$ws.Range("E12").Value = -0.000792707094728673

$ws.Range("E13").Value = -0.004083840607741807

$ws.Protect()
